$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task description text for the payment processing row (C28)
$ws.Range("C28").Value = "Xử lý thanh toán phòng đặt phòng"

# Clear the duplicate member name in D25 and merge it with D24
# (D24 "Trần Thị Ngọc Hân" now spans rows 24-25)
$ws.Range("D25").Value = ""
$ws.Range("D24:D25").Merge()

# Adjust row 29 height to match the now-shorter wrapped content
$ws.Rows.Item(29).RowHeight = 39.75
